# Apply "Corr/total marks" change to the marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# B11: correct answers count 3 -> 5
$ws.Range("B11").Value = 5

# B12: total marks 81 -> 135
$ws.Range("B12").Value = 135

# E12: "corrected/total" display text 81/84 -> 135/140
$ws.Range("E12").Value = "135/140"
